$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pH control")
$ws.Select()
